$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 11 ("NMS 2.0 Diagnostics") - append text to the "See 2 Demo Video"
# line so it also references the new 2.5 Smart Box demo video.
# ---------------------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$shp11 = $s11.Shapes.Item(2)
$tr11 = $shp11.TextFrame.TextRange
$paras11 = @($tr11.Paragraphs())
$runs11 = @($paras11[5].Runs())
$runs11[0].Text = "See 2 Demo Video (Diagnostics) and 2.5 Demo Video (Smart Box)"

# ---------------------------------------------------------------------------
# Slide 13 ("Shutdown NMS 2.0") - rewrite the body placeholder: new intro
# line, reworded/expanded stop instructions, and a new closing instruction
# about rebuilding the app.
# ---------------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$shp13 = $s13.Shapes.Item(2)
$tf13 = $shp13.TextFrame
$tr13 = $tf13.TextRange

$lines13 = @(
    "Depending on use case, use one of the following:",
    "Stop only python_app by selecting python terminal, then enter Ctrl+C (KeyboardInterrupt is coded to stop threads safely).",
    "Stop NMS 2.0 by entering “docker compose stop” in regular terminal/cmd, then “docker compose start” to restart.",
    "Delete NMS 2.0 container (not local files) by entering “docker compose down” in regular terminal/cmd.",
    "You can now make changes to nms2 files, then “docker compose up” to rebuild the updated app and see the changes.",
    "",
    "See 4 Demo Video (Shutdown)"
)
$tr13.Text = [string]::Join("`r", $lines13)

$paras13 = @($tr13.Paragraphs())

# Paragraph 0 ("Depending on use case...") and paragraph 4 ("You can now
# make changes...") / paragraph 5 (blank) have no bullet / numbering.
$paras13[0].ParagraphFormat.Bullet.Visible = 0
$paras13[4].ParagraphFormat.Bullet.Visible = 0
$paras13[5].ParagraphFormat.Bullet.Visible = 0

# --- Paragraph 1: split out the runs that carry the spell-check-exempt
# words (python_app, Ctrl+C, KeyboardInterrupt) as their own runs.
$par1 = $paras13[1]
$par1.Characters(11, 10).Font.Bold = 0     # python_app
$par1.Characters(63, 6).Font.Bold = 0      # Ctrl+C
$par1.Characters(71, 17).Font.Bold = 0     # KeyboardInterrupt

# --- Paragraph 2: split out "cmd" as its own run.
$par2 = $paras13[2]
$par2.Characters(68, 3).Font.Bold = 0      # cmd
